$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Priorisierung")
$ws3 = $wb.Worksheets.Item("Tasks")

# --- Tasks sheet: new "UserStory 6" section (email Bestaetigung) ---
# Header row 38, mirrors the existing section headers (row 34 etc.):
# copy A8:B8 from Priorisierung (already holds the right shared strings/styles)
$ws2.Range("A8:B8").Copy()
$ws3.Range("A38").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws3.Rows.Item(38).RowHeight = 30

# New task rows for this user story (column B only)
$ws3.Range("B40").Value = "Klasse Booking erstellen"
$ws3.Range("B41").Value = "Interface IBookingReposittory erstellen"
$ws3.Range("B42").Value = "Klasse BookingRepository DB erstellen"
$ws3.Range("B43").Value = "Booking Methode im BookingController erstellen"
$ws3.Range("B44").Value = "Booking View erstellen"
$ws3.Range("B45").Value = "email Bestätigung einbauen"
$ws3.Range("B39").Value = "Tabelle bookings und rooms erstellen"

# --- Selection/view state updates recorded in the diff ---
$ws2.Activate()
$ws2.Range("C8").Select()

$ws3.Activate()
$ws3.Range("B44").Select()
